$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 1035.7059
$ws.Range("I80").Value = 1312.125
$ws.Range("J80").Value = 790
$ws.Range("K80").Value = 3936.375
$ws.Range("L80").Value = 2370
$ws.Range("M80").Value = -2938.375
$ws.Range("N80").Value = -4366
$ws.Range("H83").Value = 1035.7059
$ws.Range("I83").Value = 1312.125
$ws.Range("J83").Value = 790
$ws.Range("K83").Value = 11809.125
$ws.Range("L83").Value = 7110
$ws.Range("M83").Value = -6817.125
$ws.Range("N83").Value = -17094
$ws.Range("H86").Value = 1645.3077
$ws.Range("I86").Value = 1115.75
$ws.Range("J86").Value = 8000
$ws.Range("K86").Value = 1115.75
$ws.Range("L86").Value = 8000
$ws.Range("M86").Value = 7.25
$ws.Range("N86").Value = -10246
$ws.Range("H89").Value = 1645.3077
$ws.Range("I89").Value = 1115.75
$ws.Range("J89").Value = 8000
$ws.Range("K89").Value = 5578.75
$ws.Range("L89").Value = 40000
$ws.Range("M89").Value = 37.25
$ws.Range("N89").Value = -51232
$ws.Range("H98").Value = 4922.222
$ws.Range("I98").Value = 4850
$ws.Range("J98").Value = 5500
$ws.Range("K98").Value = 4850
$ws.Range("L98").Value = 5500
$ws.Range("M98").Value = -3352
$ws.Range("N98").Value = -8496
$ws.Range("H103").Value = 1876.25
$ws.Range("J103").Value = 2752.5
$ws.Range("L103").Value = 8257.5
$ws.Range("N103").Value = -9429.5
$ws.Range("H112").Value = 5575.294
$ws.Range("J112").Value = 1592
$ws.Range("L112").Value = 4776
$ws.Range("N112").Value = -6992
$ws.Range("H122").Value = 4922.222
$ws.Range("I122").Value = 4850
$ws.Range("J122").Value = 5500
$ws.Range("K122").Value = 14550
$ws.Range("L122").Value = 16500
$ws.Range("M122").Value = -12100
$ws.Range("N122").Value = -21400
$ws.Range("H125").Value = 771.1429000000001
$ws.Range("J125").Value = 771.1429000000001
$ws.Range("L125").Value = 6940.2861
$ws.Range("N125").Value = -11860.2861
$ws.Range("H137").Value = 836009.4399999999
$ws.Range("I137").Value = 2677.3
$ws.Range("J137").Value = 1669341.6
$ws.Range("K137").Value = 8031.900000000001
$ws.Range("L137").Value = 5008024.800000001
$ws.Range("M137").Value = -5481.900000000001
$ws.Range("N137").Value = -5013124.800000001
$ws.Range("H138").Value = 3606.087
$ws.Range("I138").Value = 3410.5557
$ws.Range("J138").Value = 3653.6487
$ws.Range("K138").Value = 10231.6671
$ws.Range("L138").Value = 10960.9461
$ws.Range("M138").Value = -5091.667099999999
$ws.Range("N138").Value = -21240.9461

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 5307.4116
$ws.Range("I88").Value = 11781.2
$ws.Range("J88").Value = 2610
$ws.Range("K88").Value = 11781.2
$ws.Range("L88").Value = 2610
$ws.Range("M88").Value = -11375.2
$ws.Range("N88").Value = -3422
$ws.Range("H91").Value = 5307.4116
$ws.Range("I91").Value = 11781.2
$ws.Range("J91").Value = 2610
$ws.Range("K91").Value = 11781.2
$ws.Range("L91").Value = 2610
$ws.Range("M91").Value = -10377.2
$ws.Range("N91").Value = -5418

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H13").Value = 15453
$ws.Range("I13").Value = 10647
$ws.Range("J13").Value = 17856
$ws.Range("K13").Value = 10647
$ws.Range("L13").Value = 17856
$ws.Range("M13").Value = -10479
$ws.Range("N13").Value = -18192
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 298.57144
$ws.Range("I22").Value = 297.6
$ws.Range("J22").Value = 301
$ws.Range("K22").Value = 297.6
$ws.Range("L22").Value = 301
$ws.Range("M22").Value = 52.39999999999998
$ws.Range("N22").Value = -1001
$ws.Range("H99").Value = 3400
$ws.Range("J99").Value = 3400
$ws.Range("L99").Value = 3400
$ws.Range("N99").Value = -6396
$ws.Range("H126").Value = 3400
$ws.Range("J126").Value = 3400
$ws.Range("L126").Value = 10200
$ws.Range("N126").Value = -15140

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2594.606
$ws.Range("I68").Value = 816.36365
$ws.Range("J68").Value = 4372.8486
$ws.Range("K68").Value = 2449.09095
$ws.Range("L68").Value = 13118.5458
$ws.Range("M68").Value = -1638.09095
$ws.Range("N68").Value = -14740.5458
$ws.Range("H71").Value = 2594.606
$ws.Range("I71").Value = 816.36365
$ws.Range("J71").Value = 4372.8486
$ws.Range("K71").Value = 7347.27285
$ws.Range("L71").Value = 39355.6374
$ws.Range("M71").Value = -3291.27285
$ws.Range("N71").Value = -47467.6374
$ws.Range("H109").Value = 3364.5557
$ws.Range("I109").Value = 2750.3333
$ws.Range("K109").Value = 8250.999899999999
$ws.Range("M109").Value = -7210.999899999999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H117").Value = 62600
$ws.Range("J117").Value = 62600
$ws.Range("L117").Value = 62600
$ws.Range("N117").Value = -71778
$ws.Range("H129").Value = 82266.336
$ws.Range("I129").Value = 81999
$ws.Range("J129").Value = 82400
$ws.Range("K129").Value = 81999
$ws.Range("L129").Value = 82400
$ws.Range("M129").Value = -76999
$ws.Range("N129").Value = -92400

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 29716
$ws.Range("I63").Value = 18250
$ws.Range("J63").Value = 35449
$ws.Range("K63").Value = 18250
$ws.Range("L63").Value = 35449
$ws.Range("M63").Value = -17626
$ws.Range("N63").Value = -36697
$ws.Range("H66").Value = 29716
$ws.Range("I66").Value = 18250
$ws.Range("J66").Value = 35449
$ws.Range("K66").Value = 54750
$ws.Range("L66").Value = 106347
$ws.Range("M66").Value = -51630
$ws.Range("N66").Value = -112587
$ws.Range("H81").Value = 11767673
$ws.Range("I81").Value = 2228.5
$ws.Range("K81").Value = 4457
$ws.Range("M81").Value = -3396
$ws.Range("H84").Value = 11767673
$ws.Range("I84").Value = 2228.5
$ws.Range("K84").Value = 22285
$ws.Range("M84").Value = -16981
